$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.126417
$ws.Range("H2").Value = 0.379251
$ws.Range("I2").Value = 0.733230478333749
$ws.Range("J2").Value = 0.8047954435010685
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.191633
$ws.Range("N2").Value = 0.383266
$ws.Range("O2").Value = 0.2592126846256717
$ws.Range("P2").Value = 0.1891519879500531
$ws.Range("Q2").Value = 0.024225668961
$ws.Range("R2").Value = 0.145354013766
$ws.Range("S2").Value = 0.1900626407382565
$ws.Range("T2").Value = 0.1522286580313718

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr7"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.126417
$ws.Range("H3").Value = 0.379251
$ws.Range("I3").Value = 0.733230478333749
$ws.Range("J3").Value = 0.8047954435010685
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1040336666666667
$ws.Range("N3").Value = 0.312101
$ws.Range("O3").Value = 0.140721305976105
$ws.Range("P3").Value = 0.1540301633622589
$ws.Range("Q3").Value = 0.013151624039
$ws.Range("R3").Value = 0.118364616351
$ws.Range("S3").Value = 0.1031811504926093
$ws.Range("T3").Value = 0.1239627736356712

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr7"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.126417
$ws.Range("H4").Value = 0.379251
$ws.Range("I4").Value = 0.733230478333749
$ws.Range("J4").Value = 0.8047954435010685
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.2501216666666666
$ws.Range("N4").Value = 0.7503649999999999
$ws.Range("O4").Value = 0.3383274733460001
$ws.Range("P4").Value = 0.3703251304267574
$ws.Range("Q4").Value = 0.03161963073499999
$ws.Range("R4").Value = 0.284576676615
$ws.Range("S4").Value = 0.2480720151149364
$ws.Range("T4").Value = 0.2980359775813933

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr7"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.126417
$ws.Range("H5").Value = 0.379251
$ws.Range("I5").Value = 0.733230478333749
$ws.Range("J5").Value = 0.8047954435010685
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1935003333333334
$ws.Range("N5").Value = 0.580501
$ws.Range("O5").Value = 0.2617385360522232
$ws.Range("P5").Value = 0.2864927182609305
$ws.Range("Q5").Value = 0.024461731639
$ws.Range("R5").Value = 0.220155584751
$ws.Range("S5").Value = 0.1919146719879468
$ws.Range("T5").Value = 0.2305680342526323

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ccl21b"
$ws.Range("C6").Value = "Ccr7"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.045994
$ws.Range("H6").Value = 0.091988
$ws.Range("I6").Value = 0.266769521666251
$ws.Range("J6").Value = 0.1952045564989315
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.191633
$ws.Range("N6").Value = 0.383266
$ws.Range("O6").Value = 0.2592126846256717
$ws.Range("P6").Value = 0.1891519879500531
$ws.Range("Q6").Value = 0.008813968202
$ws.Range("R6").Value = 0.035255872808
$ws.Range("S6").Value = 0.06915004388741522
$ws.Range("T6").Value = 0.03692332991868136

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ccl21b"
$ws.Range("C7").Value = "Ccr7"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.045994
$ws.Range("H7").Value = 0.091988
$ws.Range("I7").Value = 0.266769521666251
$ws.Range("J7").Value = 0.1952045564989315
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1040336666666667
$ws.Range("N7").Value = 0.312101
$ws.Range("O7").Value = 0.140721305976105
$ws.Range("P7").Value = 0.1540301633622589
$ws.Range("Q7").Value = 0.004784924464666667
$ws.Range("R7").Value = 0.028709546788
$ws.Range("S7").Value = 0.03754015548349567
$ws.Range("T7").Value = 0.03006738972658773

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ccl21b"
$ws.Range("C8").Value = "Ccr7"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.045994
$ws.Range("H8").Value = 0.091988
$ws.Range("I8").Value = 0.266769521666251
$ws.Range("J8").Value = 0.1952045564989315
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.2501216666666666
$ws.Range("N8").Value = 0.7503649999999999
$ws.Range("O8").Value = 0.3383274733460001
$ws.Range("P8").Value = 0.3703251304267574
$ws.Range("Q8").Value = 0.01150409593666667
$ws.Range("R8").Value = 0.06902457562
$ws.Range("S8").Value = 0.09025545823106373
$ws.Range("T8").Value = 0.07228915284536416

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ccl21b"
$ws.Range("C9").Value = "Ccr7"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.045994
$ws.Range("H9").Value = 0.091988
$ws.Range("I9").Value = 0.266769521666251
$ws.Range("J9").Value = 0.1952045564989315
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1935003333333334
$ws.Range("N9").Value = 0.580501
$ws.Range("O9").Value = 0.2617385360522232
$ws.Range("P9").Value = 0.2864927182609305
$ws.Range("Q9").Value = 0.008899854331333334
$ws.Range("R9").Value = 0.05339912598800001
$ws.Range("S9").Value = 0.06982386406427638
$ws.Range("T9").Value = 0.05592468400829829

